$wb = $excel.ActiveWorkbook

# Sheet ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 408
$ws.Range("I2").Value = 408
$ws.Range("K2").Value = 408
$ws.Range("M2").Value = -295

# Sheet ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 454.5
$ws.Range("I8").Value = 143.375
$ws.Range("K8").Value = 430.125
$ws.Range("M8").Value = -291.125

# Sheet ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 960
$ws.Range("I19").Value = 900
$ws.Range("J19").Value = 975
$ws.Range("K19").Value = 900
$ws.Range("L19").Value = 975
$ws.Range("M19").Value = -725
$ws.Range("N19").Value = -1325

# Sheet ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 549.75
$ws.Range("I28").Value = 549.75
$ws.Range("K28").Value = 549.75
$ws.Range("M28").Value = -64.75

# Sheet ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("L88").ClearContents()

# Sheet ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").ClearContents()

# Sheet ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5558.6665
$ws.Range("I100").Value = 5380.4
$ws.Range("J100").Value = 6450
$ws.Range("K100").Value = 5380.4
$ws.Range("L100").Value = 6450
$ws.Range("M100").Value = -4839.4
$ws.Range("N100").Value = -7532

# Sheet ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 933
$ws.Range("I103").Value = 933
$ws.Range("K103").Value = 2799
$ws.Range("M103").Value = -2213

# Sheet ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1149.5
$ws.Range("I107").Value = 478.33334
$ws.Range("K107").Value = 478.33334
$ws.Range("M107").Value = 1441.66666

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3309.8
$ws.Range("J137").Value = 3387.25
$ws.Range("L137").Value = 10161.75
$ws.Range("N137").Value = -15261.75

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1325
$ws.Range("I45").Value = 1325
$ws.Range("K45").Value = 1325
$ws.Range("M45").Value = -948

# Sheet ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1507
$ws.Range("I97").Value = 1258.4
$ws.Range("K97").Value = 1258.4
$ws.Range("M97").Value = -762.4000000000001

# Sheet ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2803.3333
$ws.Range("I102").Value = 2803.3333
$ws.Range("K102").Value = 2803.3333
$ws.Range("M102").Value = -1181.3333

# Sheet ARM row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 92000
$ws.Range("J130").Value = 92000
$ws.Range("L130").Value = 92000
$ws.Range("N130").Value = -102040

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2812.75
$ws.Range("I132").Value = 2417
$ws.Range("K132").Value = 7251
$ws.Range("M132").Value = -4721

# Sheet BSM row 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 8499.333000000001
$ws.Range("J12").Value = 7749.5
$ws.Range("L12").Value = 7749.5
$ws.Range("N12").Value = -8085.5

# Sheet CRP row 37
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 9000
$ws.Range("J37").Value = 9000
$ws.Range("L37").Value = 9000
$ws.Range("M37").Value = -9214

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3734.875
$ws.Range("I58").Value = 3996.5
$ws.Range("J58").Value = 2950
$ws.Range("K58").Value = 3996.5
$ws.Range("L58").Value = 2950
$ws.Range("M58").Value = -3793.5
$ws.Range("N58").Value = -3356

# Sheet CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Sheet CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3500
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3500
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1132.875
$ws.Range("I134").Value = 1009
$ws.Range("K134").Value = 3027
$ws.Range("M134").Value = -492

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3734.875
$ws.Range("I136").Value = 3996.5
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 11989.5
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -9439.5
$ws.Range("N136").Value = -13950

# Sheet CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.590908
$ws.Range("I2").Value = 8.75
$ws.Range("J2").Value = 61.2
$ws.Range("K2").Value = 52.5
$ws.Range("L2").Value = 367.2
$ws.Range("M2").Value = 60.5
$ws.Range("N2").Value = -593.2

# Sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1276.6666
$ws.Range("I68").Value = 798.5
$ws.Range("K68").Value = 2395.5
$ws.Range("M68").Value = -1584.5

# Sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1276.6666
$ws.Range("I71").Value = 798.5
$ws.Range("K71").Value = 7186.5
$ws.Range("M71").Value = -3130.5

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2329.7144
$ws.Range("I131").Value = 1320
$ws.Range("J131").Value = 2498
$ws.Range("K131").Value = 3960
$ws.Range("L131").Value = 7494
$ws.Range("M131").Value = 1080
$ws.Range("N131").Value = -17574

# Sheet CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 816.6667
$ws.Range("I140").Value = 825
$ws.Range("J140").Value = 800
$ws.Range("K140").Value = 2475
$ws.Range("L140").Value = 2400
$ws.Range("M140").Value = 2705
$ws.Range("N140").Value = -12760

# Sheet GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.63636
$ws.Range("I2").Value = 111.42857
$ws.Range("K2").Value = 111.42857
$ws.Range("M2").Value = 1.571430000000007

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 939.8
$ws.Range("I97").Value = 939.8
$ws.Range("K97").Value = 939.8
$ws.Range("M97").Value = -443.8

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2067.087
$ws.Range("I102").Value = 2128.158
$ws.Range("K102").Value = 2128.158
$ws.Range("M102").Value = -506.1579999999999

# Sheet GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("L113").ClearContents()

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1118.8889
$ws.Range("I122").Value = 1118.8889
$ws.Range("K122").Value = 3356.6667
$ws.Range("M122").Value = -906.6666999999998

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4590.143
$ws.Range("I132").Value = 4174
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 12522
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -9992
$ws.Range("N132").Value = -35060

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7599.6665
$ws.Range("I40").Value = 7599.6665
$ws.Range("K40").Value = 7599.6665
$ws.Range("M40").Value = -7463.6665

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4167.683
$ws.Range("I46").Value = 2530.0833
$ws.Range("J46").Value = 4845.3105
$ws.Range("K46").Value = 2530.0833
$ws.Range("L46").Value = 4845.3105
$ws.Range("M46").Value = -2342.0833
$ws.Range("N46").Value = -5221.3105

# Sheet LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 888.5
$ws.Range("I93").Value = 400
$ws.Range("J93").Value = 1377
$ws.Range("K93").Value = 400
$ws.Range("L93").Value = 1377
$ws.Range("M93").Value = 848
$ws.Range("N93").Value = -3873

# Sheet LTW row 99
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# Sheet LTW row 104
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 30001
$ws.Range("J104").Value = 30001
$ws.Range("L104").Value = 30001
$ws.Range("N104").Value = -36989

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6079.4
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# Sheet WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 89796.5
$ws.Range("J135").Value = 89796.5
$ws.Range("L135").Value = 89796.5
$ws.Range("N135").Value = -99936.5
